$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update aggregate metrics to reflect the newly closed trade
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.67   # Current Capital
$summary.Range("B4").Value = 0.66      # Total P&L $
$summary.Range("B5").Value = 0.49      # Total P&L %
$summary.Range("B6").Value = 27        # Total Trades
$summary.Range("B8").Value = 8         # Losing Trades
$summary.Range("B9").Value = 40.74     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.67     # Capital
$status.Range("D4").Value = 27         # Trades
$status.Range("E4").Value = 0.66       # P&L $
$status.Range("F4").Value = 0.67       # P&L %
$status.Range("G4").Value = 40.74      # Win Rate %

# ---------------------------------------------------------------------------
# Append the newly-closed trade (#27) as row 28 on both the "All Trades" and
# "MarketMaking" sheets (they mirror each other).
# ---------------------------------------------------------------------------
function Add-TradeRow($ws) {
    $ws.Cells.Item(28, 1).Value = 27

    # Columns B/C hold date/time-looking text that must stay plain text
    # (matching the rest of the sheet) rather than being auto-converted to
    # Excel date/time serials.
    $ws.Cells.Item(28, 2).NumberFormat = "@"
    $ws.Cells.Item(28, 2).Value = "2026-02-17"
    $ws.Cells.Item(28, 3).NumberFormat = "@"
    $ws.Cells.Item(28, 3).Value = "12:37:45"

    $ws.Cells.Item(28, 4).Value = "MarketMaking"
    $ws.Cells.Item(28, 5).Value = "DOWN"
    $ws.Cells.Item(28, 6).Value = 0.39
    $ws.Cells.Item(28, 7).Value = 0.35
    $ws.Cells.Item(28, 8).Value = "CLOSED"
    $ws.Cells.Item(28, 9).Value = -10.2564
    $ws.Cells.Item(28, 10).Value = -0.04
    $ws.Cells.Item(28, 11).Value = 100.67
    $ws.Cells.Item(28, 12).Value = 0
    $ws.Cells.Item(28, 13).Value = 0
    $ws.Cells.Item(28, 14).Value = 0.6
    $ws.Cells.Item(28, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(28, 16).Value = "early_exit"
    $ws.Cells.Item(28, 17).Value = 0.14
}

Add-TradeRow $wb.Worksheets.Item("All Trades")
Add-TradeRow $wb.Worksheets.Item("MarketMaking")
